$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 9

# Update values in row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3

# Remove rows 5 and 6 entirely (shrinking the used range to A1:B4)
$ws.Range("A5:B6").Delete()
